$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2007-2009 rows (rows 2-4); remaining data shifts up.
$ws.Rows("2:4").Delete()

# Append the new 2021 row of data as row 13.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 121
$ws.Range("F13").Value = 245
$ws.Range("H13").Value = 51
$ws.Range("K13").Value = 30
$ws.Range("L13").Value = 1
